$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '42.025.37'
Set-TextValue $ws.Range('E2') '  +0.53%  '
Set-TextValue $ws.Range('D3') '2.264.79'
Set-TextValue $ws.Range('E3') '  -0.31%  '
Set-TextValue $ws.Range('E4') '  -0.01%  '
Set-TextValue $ws.Range('D5') '153.11'
Set-TextValue $ws.Range('E5') '  +15,194.06%  '
Set-TextValue $ws.Range('D6') '305.36'
Set-TextValue $ws.Range('E6') '  +0.97%  '
Set-TextValue $ws.Range('D7') '93.95'
Set-TextValue $ws.Range('E7') '  +1.83%  '
Set-TextValue $ws.Range('D8') '0.529'
Set-TextValue $ws.Range('E8') '  -0.35%  '
Set-TextValue $ws.Range('E9') '  +0.00%  '
Set-TextValue $ws.Range('E10') '  +0.08%  '
Set-TextValue $ws.Range('D11') '33.57'
Set-TextValue $ws.Range('E11') '  +4.17%  '
Set-TextValue $ws.Range('D12') '0.0799'
Set-TextValue $ws.Range('E12') '  -0.18%  '
Set-TextValue $ws.Range('E13') '  -2.32%  '
Set-TextValue $ws.Range('D14') '6.66'
Set-TextValue $ws.Range('E14') '  -0.23%  '
Set-TextValue $ws.Range('D15') '2.616.59'
Set-TextValue $ws.Range('E15') '  -0.23%  '
Set-TextValue $ws.Range('D16') '14.30'
Set-TextValue $ws.Range('E16') '  +0.59%  '
Set-TextValue $ws.Range('D17') '2.269.12'
Set-TextValue $ws.Range('E18') '  +3.70%  '
Set-TextValue $ws.Range('D19') '41.956.30'
Set-TextValue $ws.Range('E19') '  +0.56%  '
Set-TextValue $ws.Range('D20') '12.86'
Set-TextValue $ws.Range('E20') '  +5.70%  '
Set-TextValue $ws.Range('E21') '  +0.92%  '
Set-TextValue $ws.Range('D22') '5.98'
Set-TextValue $ws.Range('E22') '  +0.42%  '
Set-TextValue $ws.Range('D23') '67.93'
Set-TextValue $ws.Range('E23') '  +1.16%  '
Set-TextValue $ws.Range('D24') '243.07'
Set-TextValue $ws.Range('E24') '  -0.04%  '
Set-TextValue $ws.Range('E25') '  +1.08%  '
Set-TextValue $ws.Range('D26') '1.93'
Set-TextValue $ws.Range('E26') '  +1.23%  '
Set-TextValue $ws.Range('E27') '  -0.02%  '
Set-TextValue $ws.Range('D28') '23.95'
Set-TextValue $ws.Range('E28') '  -1.13%  '
Set-TextValue $ws.Range('D29') '9.71'
Set-TextValue $ws.Range('E29') '  +0.83%  '
Set-TextValue $ws.Range('B30') 'InjectiveProtocol'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D30') '35.40'
Set-TextValue $ws.Range('E30') '  +4.79%  '
Set-TextValue $ws.Range('B31') 'Toncoin'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D31') '2.09'
Set-TextValue $ws.Range('E31') '  -0.05%  '
Set-TextValue $ws.Range('D32') '159.83'
Set-TextValue $ws.Range('E32') '  +0.48%  '
Set-TextValue $ws.Range('E33') '  +2.91%  '
Set-TextValue $ws.Range('D34') '1.00'
Set-TextValue $ws.Range('E34') '  +0.01%  '
Set-TextValue $ws.Range('D35') '0.0748'
Set-TextValue $ws.Range('E35') '  +0.16%  '
Set-TextValue $ws.Range('D36') '3.08'
Set-TextValue $ws.Range('E36') '  +0.38%  '
Set-TextValue $ws.Range('D37') '17.11'
Set-TextValue $ws.Range('E37') '  +2.15%  '
Set-TextValue $ws.Range('D38') '0.107'
Set-TextValue $ws.Range('E38') '  +2.62%  '
Set-TextValue $ws.Range('E39') '  -0.83%  '
Set-TextValue $ws.Range('E40') '  +0.18%  '
Set-TextValue $ws.Range('E41') '  -1.41%  '
Set-TextValue $ws.Range('D42') '4.12'
Set-TextValue $ws.Range('E42') '  +5.00%  '
Set-TextValue $ws.Range('D43') '19.73'
Set-TextValue $ws.Range('E43') '  +0.39%  '
Set-TextValue $ws.Range('D44') '2.001.50'
Set-TextValue $ws.Range('E44') '  -3.70%  '
Set-TextValue $ws.Range('E45') '  +10.79%  '
Set-TextValue $ws.Range('E46') '  +1.13%  '
Set-TextValue $ws.Range('D47') '10.22'
Set-TextValue $ws.Range('E47') '  +0.15%  '
Set-TextValue $ws.Range('D48') '2.92'
Set-TextValue $ws.Range('E48') '  -0.91%  '
Set-TextValue $ws.Range('D49') '53.51'
Set-TextValue $ws.Range('E49') '  +2.37%  '
Set-TextValue $ws.Range('D50') '72.56'
Set-TextValue $ws.Range('E50') '  -1.63%  '
Set-TextValue $ws.Range('D51') '1.51'
Set-TextValue $ws.Range('E51') '  -0.99%  '
